# Fri, May 01, 2020  2:06:37 PM
#
# 1) Three tables (slides 14, 15, 16) get their table style switched from
#    the custom "Table_0" style to the built-in style
#    {44A245FD-2DEC-493E-8607-72E1929F3A12}.
# 2) The presentation's theme colour scheme ("Integral" / "Red Violet")
#    is swapped for the classic "Office" colour scheme - i.e. the colours
#    that used to live in the (otherwise unused) secondary theme part end
#    up driving the deck.

$p = $ppt.ActivePresentation

# --- 1. Update the table styles on the three tables --------------------
$newStyleId = "{44A245FD-2DEC-493E-8607-72E1929F3A12}"
$tableSlideIndexes = @(14, 15, 16)

foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Swap the theme colour scheme to the "Office" palette -----------
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Office theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM RGB() values (R + G*256 + B*65536).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
